$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of search-result data (row 4), matching the format of
# existing rows 2-3.
$row = 4

$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item($row, 1).Value = 42602.583657407406

$ws.Cells.Item($row, 2).Value = "Noun"

$values = @(12267, 6887, 1277, 156, 73, 67, 31, 4, 3, 57, 42)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($row, 3 + $i).Value = $values[$i]
}
